$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 111706580
$ws.Range("B2").Value = 88032
$ws.Range("D2").Value = "VU"
$ws.Range("E2").Value = 6276
$ws.Range("F2").Value = "Goliatmusseron"
$ws.Range("G2").Value = "Tricholoma matsutake"
$ws.Range("H2").Value = "(S.Ito & S.Imai) Singer"
$ws.Range("Q2").Value = 754511
$ws.Range("R2").Value = 7089952
$ws.Range("Z2").Value = "14:48"
$ws.Range("AB2").Value = "14:48"

# Row 3
$ws.Range("A3").Value = 111708888
$ws.Range("B3").Value = 90678
$ws.Range("D3").Value = "LC"
$ws.Range("E3").Value = 4366
$ws.Range("F3").Value = "Skarp dropptaggsvamp"
$ws.Range("G3").Value = "Hydnellum peckii"
$ws.Range("H3").Value = "Banker"
$ws.Range("Q3").Value = 754511
$ws.Range("R3").Value = 7089952
$ws.Range("Z3").Value = "13:54"
$ws.Range("AB3").Value = "13:54"

# Row 4
$ws.Range("A4").Value = 111708099
$ws.Range("B4").Value = 90660
$ws.Range("D4").Value = "NT"
$ws.Range("E4").Value = 4362
$ws.Range("F4").Value = "Blå taggsvamp"
$ws.Range("G4").Value = "Hydnellum caeruleum"
$ws.Range("H4").Value = "(Hornem.) P.Karst."
$ws.Range("Q4").Value = 754511
$ws.Range("R4").Value = 7089952
$ws.Range("Z4").Value = "14:16"
$ws.Range("AB4").Value = "14:16"

# Row 5
$ws.Range("A5").Value = 111708040
$ws.Range("Q5").Value = 754511
$ws.Range("R5").Value = 7089952
$ws.Range("Z5").Value = "14:20"
$ws.Range("AB5").Value = "14:20"

# Row 6
$ws.Range("A6").Value = 111708029
$ws.Range("B6").Value = 90662
$ws.Range("D6").Value = "LC"
$ws.Range("E6").Value = 4363
$ws.Range("F6").Value = "Zontaggsvamp"
$ws.Range("G6").Value = "Hydnellum concrescens"
$ws.Range("H6").Value = "(Pers.) Banker"
$ws.Range("Q6").Value = 754511
$ws.Range("R6").Value = 7089952
$ws.Range("Z6").Value = "14:21"
$ws.Range("AB6").Value = "14:21"

# Row 7
$ws.Range("A7").Value = 111704319
$ws.Range("B7").Value = 90710
$ws.Range("D7").Value = "NT"
$ws.Range("E7").Value = 5449
$ws.Range("F7").Value = "Svart taggsvamp"
$ws.Range("G7").Value = "Phellodon niger"
$ws.Range("H7").Value = "(Fr.:Fr.) P.Karst."
$ws.Range("Q7").Value = 754511
$ws.Range("R7").Value = 7089952
$ws.Range("Z7").Value = "15:11"
$ws.Range("AB7").Value = "15:11"

# Row 8
$ws.Range("A8").Value = 111708126
$ws.Range("B8").Value = 90666
$ws.Range("D8").Value = "LC"
$ws.Range("E8").Value = 4364
$ws.Range("F8").Value = "Dropptaggsvamp"
$ws.Range("G8").Value = "Hydnellum ferrugineum"
$ws.Range("H8").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("Q8").Value = 754511
$ws.Range("R8").Value = 7089952
$ws.Range("Z8").Value = "14:14"
$ws.Range("AB8").Value = "14:14"

# Row 9
$ws.Range("Q9").Value = 754511
$ws.Range("R9").Value = 7089952
